$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 5134.5
$ws.Cells.Item(62, 9).Value = 3040
$ws.Cells.Item(62, 10).Value = 8625.333000000001
$ws.Cells.Item(62, 11).Value = 3040
$ws.Cells.Item(62, 12).Value = 8625.333000000001
$ws.Cells.Item(62, 13).Value = -2416
$ws.Cells.Item(62, 14).Value = -9873.333000000001
$ws.Cells.Item(65, 8).Value = 5134.5
$ws.Cells.Item(65, 9).Value = 3040
$ws.Cells.Item(65, 10).Value = 8625.333000000001
$ws.Cells.Item(65, 11).Value = 15200
$ws.Cells.Item(65, 12).Value = 43126.665
$ws.Cells.Item(65, 13).Value = -12080
$ws.Cells.Item(65, 14).Value = -49366.665
$ws.Cells.Item(100, 8).Value = 25001882
$ws.Cells.Item(100, 9).Value = 28573422
$ws.Cells.Item(100, 11).Value = 28573422
$ws.Cells.Item(100, 13).Value = -28572881
$ws.Cells.Item(106, 8).Value = 27780824
$ws.Cells.Item(106, 9).Value = 1420.8334
$ws.Cells.Item(106, 11).Value = 1420.8334
$ws.Cells.Item(106, 13).Value = -789.8334
$ws.Cells.Item(129, 8).Value = 859.567
$ws.Cells.Item(129, 10).Value = 891.8022
$ws.Cells.Item(129, 12).Value = 2675.4066
$ws.Cells.Item(129, 14).Value = -12675.4066
$ws.Cells.Item(132, 8).Value = 335104.5
$ws.Cells.Item(132, 9).Value = 1341.8148
$ws.Cells.Item(132, 10).Value = 3338968.8
$ws.Cells.Item(132, 11).Value = 4025.4444
$ws.Cells.Item(132, 12).Value = 10016906.4
$ws.Cells.Item(132, 13).Value = -1495.4444
$ws.Cells.Item(132, 14).Value = -10021966.4
$ws.Cells.Item(138, 8).Value = 5454.73
$ws.Cells.Item(138, 10).Value = 7120.5137
$ws.Cells.Item(138, 12).Value = 21361.5411
$ws.Cells.Item(138, 14).Value = -31641.5411

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 798.2727
$ws.Cells.Item(2, 9).Value = 651.8823
$ws.Cells.Item(2, 10).Value = 1296
$ws.Cells.Item(2, 11).Value = 651.8823
$ws.Cells.Item(2, 12).Value = 1296
$ws.Cells.Item(2, 13).Value = -538.8823
$ws.Cells.Item(2, 14).Value = -1522
$ws.Cells.Item(6, 8).Value = 14829.333
$ws.Cells.Item(6, 10).Value = 14829.333
$ws.Cells.Item(6, 12).Value = 14829.333
$ws.Cells.Item(6, 14).Value = -15175.333
$ws.Cells.Item(45, 8).Value = 1983
$ws.Cells.Item(45, 9).Value = 2343.8333
$ws.Cells.Item(45, 10).Value = 1622.1666
$ws.Cells.Item(45, 11).Value = 2343.8333
$ws.Cells.Item(45, 12).Value = 1622.1666
$ws.Cells.Item(45, 13).Value = -1966.8333
$ws.Cells.Item(45, 14).Value = -2376.1666
$ws.Cells.Item(61, 8).Value = 1689.5
$ws.Cells.Item(61, 9).Value = 1431.8
$ws.Cells.Item(61, 10).Value = 2978
$ws.Cells.Item(61, 11).Value = 1431.8
$ws.Cells.Item(61, 12).Value = 2978
$ws.Cells.Item(61, 13).Value = -1219.8
$ws.Cells.Item(61, 14).Value = -3402
$ws.Cells.Item(116, 8).Value = 798.2727
$ws.Cells.Item(116, 9).Value = 651.8823
$ws.Cells.Item(116, 10).Value = 1296
$ws.Cells.Item(116, 11).Value = 651.8823
$ws.Cells.Item(116, 12).Value = 1296
$ws.Cells.Item(116, 13).Value = 1642.1177
$ws.Cells.Item(116, 14).Value = -5884
$ws.Cells.Item(136, 8).Value = 1689.5
$ws.Cells.Item(136, 9).Value = 1431.8
$ws.Cells.Item(136, 10).Value = 2978
$ws.Cells.Item(136, 11).Value = 4295.4
$ws.Cells.Item(136, 12).Value = 8934
$ws.Cells.Item(136, 13).Value = -1745.4
$ws.Cells.Item(136, 14).Value = -14034
$ws.Cells.Item(137, 8).Value = 44587.5
$ws.Cells.Item(137, 10).Value = 44587.5
$ws.Cells.Item(137, 12).Value = 44587.5
$ws.Cells.Item(137, 14).Value = -54787.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 798.2727
$ws.Cells.Item(3, 9).Value = 651.8823
$ws.Cells.Item(3, 10).Value = 1296
$ws.Cells.Item(3, 11).Value = 651.8823
$ws.Cells.Item(3, 12).Value = 1296
$ws.Cells.Item(3, 13).Value = -537.8823
$ws.Cells.Item(3, 14).Value = -1524
$ws.Cells.Item(8, 8).Value = 5093
$ws.Cells.Item(8, 9).Value = 1264.5
$ws.Cells.Item(8, 10).Value = 12750
$ws.Cells.Item(8, 11).Value = 1264.5
$ws.Cells.Item(8, 12).Value = 12750
$ws.Cells.Item(8, 13).Value = -1124.5
$ws.Cells.Item(8, 14).Value = -13030
$ws.Cells.Item(12, 8).Value = 3266.6667
$ws.Cells.Item(12, 9).Value = 650
$ws.Cells.Item(12, 10).Value = 8500
$ws.Cells.Item(12, 11).Value = 650
$ws.Cells.Item(12, 12).Value = 8500
$ws.Cells.Item(12, 13).Value = -482
$ws.Cells.Item(12, 14).Value = -8836
$ws.Cells.Item(137, 8).Value = 45370
$ws.Cells.Item(137, 10).Value = 45370
$ws.Cells.Item(137, 12).Value = 45370
$ws.Cells.Item(137, 14).Value = -55570

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 11112154
$ws.Cells.Item(16, 9).Value = 18519296
$ws.Cells.Item(16, 10).Value = 1443.75
$ws.Cells.Item(16, 11).Value = 18519296
$ws.Cells.Item(16, 12).Value = 1443.75
$ws.Cells.Item(16, 13).Value = -18519009
$ws.Cells.Item(16, 14).Value = -2017.75
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 14).ClearContents() # N37
$ws.Cells.Item(52, 8).Value = 39000
$ws.Cells.Item(52, 10).Value = 39000
$ws.Cells.Item(52, 12).Value = 39000
$ws.Cells.Item(52, 14).Value = -39588
$ws.Cells.Item(93, 8).Value = 7101.75
$ws.Cells.Item(93, 9).Value = 7101.75
$ws.Cells.Item(93, 11).Value = 7101.75
$ws.Cells.Item(93, 13).Value = -5229.75
$ws.Cells.Item(113, 8).Value = 11112154
$ws.Cells.Item(113, 9).Value = 18519296
$ws.Cells.Item(113, 10).Value = 1443.75
$ws.Cells.Item(113, 11).Value = 18519296
$ws.Cells.Item(113, 12).Value = 1443.75
$ws.Cells.Item(113, 13).Value = -18517126
$ws.Cells.Item(113, 14).Value = -5783.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 30232.25
$ws.Cells.Item(11, 9).Value = 30232.25
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 90696.75
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = -90556.75
$ws.Cells.Item(11, 14).ClearContents() # N11
$ws.Cells.Item(129, 8).Value = 2812.0908
$ws.Cells.Item(129, 9).Value = 4066.6667
$ws.Cells.Item(129, 11).Value = 12200.0001
$ws.Cells.Item(129, 13).Value = -7200.000100000001
$ws.Cells.Item(132, 8).Value = 2476.8462
$ws.Cells.Item(132, 9).Value = 944.44446
$ws.Cells.Item(132, 10).Value = 5924.75
$ws.Cells.Item(132, 11).Value = 8500.00014
$ws.Cells.Item(132, 12).Value = 53322.75
$ws.Cells.Item(132, 13).Value = -5970.00014
$ws.Cells.Item(132, 14).Value = -58382.75
$ws.Cells.Item(133, 8).Value = 2956
$ws.Cells.Item(133, 9).Value = 2495
$ws.Cells.Item(133, 11).Value = 7485
$ws.Cells.Item(133, 13).Value = -2425

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 9752817
$ws.Cells.Item(11, 10).Value = 1504695.6
$ws.Cells.Item(11, 12).Value = 1504695.6
$ws.Cells.Item(11, 14).Value = -1504973.6
$ws.Cells.Item(70, 8).Value = 5576.5366
$ws.Cells.Item(70, 9).Value = 5202.909
$ws.Cells.Item(70, 10).Value = 7117.75
$ws.Cells.Item(70, 11).Value = 5202.909
$ws.Cells.Item(70, 12).Value = 7117.75
$ws.Cells.Item(70, 13).Value = -4932.909
$ws.Cells.Item(70, 14).Value = -7657.75
$ws.Cells.Item(73, 8).Value = 5576.5366
$ws.Cells.Item(73, 9).Value = 5202.909
$ws.Cells.Item(73, 10).Value = 7117.75
$ws.Cells.Item(73, 11).Value = 5202.909
$ws.Cells.Item(73, 12).Value = 7117.75
$ws.Cells.Item(73, 13).Value = -4266.909
$ws.Cells.Item(73, 14).Value = -8989.75
$ws.Cells.Item(93, 8).Value = 25999
$ws.Cells.Item(93, 10).Value = 25999
$ws.Cells.Item(93, 12).Value = 25999
$ws.Cells.Item(93, 14).Value = -29743
$ws.Cells.Item(137, 8).Value = 38650
$ws.Cells.Item(137, 10).Value = 38650
$ws.Cells.Item(137, 12).Value = 38650
$ws.Cells.Item(137, 14).Value = -48850

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(127, 8).Value = 31470
$ws.Cells.Item(127, 10).Value = 31470
$ws.Cells.Item(127, 12).Value = 31470
$ws.Cells.Item(127, 14).Value = -41390

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 95000
$ws.Cells.Item(46, 10).Value = 95000
$ws.Cells.Item(46, 12).Value = 95000
$ws.Cells.Item(46, 14).Value = -95462
$ws.Cells.Item(96, 8).Value = 202100600
$ws.Cells.Item(96, 9).Value = 202100600
$ws.Cells.Item(96, 11).Value = 202100600
$ws.Cells.Item(96, 13).Value = -202099227
$ws.Cells.Item(113, 8).Value = 286.85715
$ws.Cells.Item(113, 9).Value = 82.166664
$ws.Cells.Item(113, 11).Value = 246.499992
$ws.Cells.Item(113, 13).Value = 1923.500008
$ws.Cells.Item(126, 8).Value = 2271.875
$ws.Cells.Item(126, 9).Value = 1380
$ws.Cells.Item(126, 10).Value = 3418.5715
$ws.Cells.Item(126, 11).Value = 4140
$ws.Cells.Item(126, 12).Value = 10255.7145
$ws.Cells.Item(126, 13).Value = -1670
$ws.Cells.Item(126, 14).Value = -15195.7145
$ws.Cells.Item(132, 8).Value = 4196.7144
$ws.Cells.Item(132, 9).Value = 2719.875
$ws.Cells.Item(132, 10).Value = 6165.8335
$ws.Cells.Item(132, 11).Value = 8159.625
$ws.Cells.Item(132, 12).Value = 18497.5005
$ws.Cells.Item(132, 13).Value = -5629.625
$ws.Cells.Item(132, 14).Value = -23557.5005
$ws.Cells.Item(134, 8).Value = 95000
$ws.Cells.Item(134, 10).Value = 95000
$ws.Cells.Item(134, 12).Value = 285000
$ws.Cells.Item(134, 14).Value = -290070
